# New code by shambhu
# Adds a "DomesticFTA" worksheet (after "Login") with a small data-binding
# table used for a domestic funds-transfer test, and selects it as the
# active sheet/tab.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login")

# Insert the new sheet right after "Login" so tab order is Login, DomesticFTA.
$ws = $wb.Worksheets.Add($null, $loginSheet)
$ws.Name = "DomesticFTA"

# --- Populate the sheet, in the same cell order the data was originally
# --- entered, so that new shared-string entries line up with the source
# --- workbook (ReciverBankName, ReceiverName, ... then the row data).
$ws.Range("B1").Value = "ReciverBankName"
$ws.Range("C1").Value = "ReceiverName"
$ws.Range("D1").Value = "ReceAccNum"
$ws.Range("E1").Value = "SwiftMsg"
$ws.Range("G1").Value = "TransferType"
$ws.Range("H1").Value = "DOT"
$ws.Range("I1").Value = "TransferDesc"

$ws.Range("B2").Value = "SBI"
$ws.Range("E2").Value = "MT103"
$ws.Range("G2").Value = "Domestic Transfer"

$ws.Range("A1").Value = "DataBinding"
$ws.Range("A2").Value = "Data001"
$ws.Range("C2").Value = "John"
$ws.Range("F1").Value = "Amount"

$ws.Range("A3").Value = "Data002"
$ws.Range("A4").Value = "Data003"
$ws.Range("A5").Value = "Data004"
$ws.Range("A6").Value = "Data005"

$ws.Range("B3").Value = "HDFC"
$ws.Range("B4").Value = "RBS"
$ws.Range("B5").Value = "ICICI"
$ws.Range("B6").Value = "CITI"

$ws.Range("C6").Value = "Pitter"
$ws.Range("C5").Value = "Raj"
$ws.Range("C4").Value = "Smith"
$ws.Range("C3").Value = "Jacson"

# Remaining repeated text cells (reuse existing shared strings).
$ws.Range("I2").Value = "TransferDesc"
$ws.Range("I3").Value = "TransferDesc"
$ws.Range("I4").Value = "TransferDesc"
$ws.Range("I5").Value = "TransferDesc"
$ws.Range("I6").Value = "TransferDesc"

$ws.Range("E3").Value = "MT103"
$ws.Range("E4").Value = "MT103"
$ws.Range("E5").Value = "MT103"
$ws.Range("E6").Value = "MT103"

$ws.Range("G3").Value = "Domestic Transfer"
$ws.Range("G4").Value = "Domestic Transfer"
$ws.Range("G5").Value = "Domestic Transfer"
$ws.Range("G6").Value = "Domestic Transfer"

# Account numbers (plain numbers).
$ws.Range("D2").Value = 1234556655
$ws.Range("D3").Value = 1234556656
$ws.Range("D4").Value = 1234556657
$ws.Range("D5").Value = 1234556658
$ws.Range("D6").Value = 1234556659

# Amounts (plain numbers).
$ws.Range("F2").Value = 8
$ws.Range("F3").Value = 10
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 9
$ws.Range("F6").Value = 2

# Date-of-transfer column: apply a date number format to H2 first, then
# copy that formatted cell down so every row shares one style record
# (rather than minting a new style per cell), then overwrite the values.
$ws.Range("H2").Value = 117924
$ws.Range("H2").NumberFormat = "mm-dd-yy"
$ws.Range("H2").Copy($ws.Range("H3:H6"))
$ws.Range("H3").Value = 118655
$ws.Range("H4").Value = 120481
$ws.Range("H5").Value = 118655
$ws.Range("H6").Value = 117559

# Approximate column widths to roughly match the authored layout.
$ws.Columns.Item(2).ColumnWidth = 17.43
$ws.Columns.Item(3).ColumnWidth = 14.14
$ws.Columns.Item(4).ColumnWidth = 14.14
$ws.Columns.Item(5).ColumnWidth = 12.57
$ws.Columns.Item(7).ColumnWidth = 17.29
$ws.Columns.Item(8).ColumnWidth = 12.86

# Selection / active state, matching the authored view.
[void]$ws.Range("K13").Select()

Write-Output "DomesticFTA sheet added"
